$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1174
$ws.Range("G3").Value = 68
$ws.Range("F7").Value = 837
$ws.Range("F8").Value = 437
$ws.Range("F10").Value = 2087
$ws.Range("F11").Value = 592
$ws.Range("F12").Value = 255
$ws.Range("F14").Value = 981
$ws.Range("F16").Value = 2113
$ws.Range("F17").Value = 587
$ws.Range("F18").Value = 10407
$ws.Range("F19").Value = 1030
$ws.Range("F20").Value = 542
$ws.Range("F21").Value = 101
$ws.Range("F26").Value = 2

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 534
$ws.Range("F10").Value = 139
$ws.Range("F11").Value = 10
$ws.Range("F13").Value = 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5656
$ws.Range("F3").Value = 457

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5656
$ws.Range("F4").Value = 457
$ws.Range("F6").Value = 534
$ws.Range("F7").Value = 1174
$ws.Range("G7").Value = 68
$ws.Range("F12").Value = 837
$ws.Range("F14").Value = 437
$ws.Range("F16").Value = 2087
$ws.Range("F17").Value = 592
$ws.Range("F18").Value = 255
$ws.Range("F22").Value = 981
$ws.Range("F25").Value = 139
$ws.Range("F26").Value = 10
$ws.Range("F27").Value = 2113
$ws.Range("F28").Value = 587
$ws.Range("F30").Value = 42
$ws.Range("F31").Value = 1030
$ws.Range("F32").Value = 542
$ws.Range("F33").Value = 101
$ws.Range("F45").Value = 2
